$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue $ws.Range("D2") "275.57"
Set-TextValue $ws.Range("G2") "22"
Set-TextValue $ws.Range("D3") "21.14"
Set-TextValue $ws.Range("G3") "22"
Set-TextValue $ws.Range("D4") "6.268"
Set-TextValue $ws.Range("G4") "22"
Set-TextValue $ws.Range("D5") "0.06235"
Set-TextValue $ws.Range("G5") "22"
Set-TextValue $ws.Range("D6") "3.556"
Set-TextValue $ws.Range("G6") "22"
Set-TextValue $ws.Range("D7") "1.531"
Set-TextValue $ws.Range("G7") "22"
Set-TextValue $ws.Range("D8") "6.564"
Set-TextValue $ws.Range("G8") "22"
Set-TextValue $ws.Range("D9") "0.8242"
Set-TextValue $ws.Range("G9") "22"
Set-TextValue $ws.Range("D10") "0.1655"
Set-TextValue $ws.Range("G10") "22"
Set-TextValue $ws.Range("D11") "0.08321"
Set-TextValue $ws.Range("G11") "22"
Set-TextValue $ws.Range("D12") "0.03497"
Set-TextValue $ws.Range("G12") "22"
Set-TextValue $ws.Range("D13") "0.03166"
Set-TextValue $ws.Range("G13") "22"
Set-TextValue $ws.Range("D14") "0.09153"
Set-TextValue $ws.Range("G14") "22"
Set-TextValue $ws.Range("D15") "3.767"
Set-TextValue $ws.Range("G15") "22"
Set-TextValue $ws.Range("D16") "0.001649"
Set-TextValue $ws.Range("G16") "22"
Set-TextValue $ws.Range("D17") "0.04667"
Set-TextValue $ws.Range("G17") "22"
Set-TextValue $ws.Range("D18") "0.006285"
Set-TextValue $ws.Range("G18") "22"
Set-TextValue $ws.Range("D19") "0.006227"
Set-TextValue $ws.Range("G19") "22"
Set-TextValue $ws.Range("G20") "22"
Set-TextValue $ws.Range("D21") "0.0001498"
Set-TextValue $ws.Range("G21") "22"
Set-TextValue $ws.Range("D22") "3.723"
Set-TextValue $ws.Range("G22") "22"
Set-TextValue $ws.Range("D23") "2.313"
Set-TextValue $ws.Range("G23") "22"
Set-TextValue $ws.Range("D24") "0.01399"
Set-TextValue $ws.Range("G24") "22"
Set-TextValue $ws.Range("G25") "22"
Set-TextValue $ws.Range("D26") "0.1242"
Set-TextValue $ws.Range("G26") "22"
Set-TextValue $ws.Range("G27") "22"
Set-TextValue $ws.Range("D28") "0.0002734"
Set-TextValue $ws.Range("G28") "22"
Set-TextValue $ws.Range("G29") "22"
Set-TextValue $ws.Range("G30") "22"
Set-TextValue $ws.Range("G31") "22"
Set-TextValue $ws.Range("G32") "22"
Set-TextValue $ws.Range("G33") "22"
Set-TextValue $ws.Range("G34") "22"
Set-TextValue $ws.Range("G35") "22"
Set-TextValue $ws.Range("G36") "22"
Set-TextValue $ws.Range("G37") "22"
Set-TextValue $ws.Range("G38") "22"
Set-TextValue $ws.Range("G39") "22"
Set-TextValue $ws.Range("D40") "0.04743"
Set-TextValue $ws.Range("G40") "22"
Set-TextValue $ws.Range("D41") "0.005293"
Set-TextValue $ws.Range("G41") "22"
Set-TextValue $ws.Range("D42") "0.007026"
Set-TextValue $ws.Range("G42") "22"
Set-TextValue $ws.Range("D43") "0.1121"
Set-TextValue $ws.Range("G43") "22"
Set-TextValue $ws.Range("D44") "0.01132"
Set-TextValue $ws.Range("G44") "22"
Set-TextValue $ws.Range("D45") "0.00006182"
Set-TextValue $ws.Range("G45") "22"
Set-TextValue $ws.Range("G46") "22"
Set-TextValue $ws.Range("D47") "0.7223"
Set-TextValue $ws.Range("G47") "22"
Set-TextValue $ws.Range("D48") "0.001395"
Set-TextValue $ws.Range("G48") "22"
Set-TextValue $ws.Range("D49") "0.00001898"
Set-TextValue $ws.Range("G49") "22"
Set-TextValue $ws.Range("G50") "22"
Set-TextValue $ws.Range("G51") "22"
